$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.020676
$ws.Range("H2").Value = 0.062028
$ws.Range("M2").Value = 4.407279333333333
$ws.Range("N2").Value = 13.221838
$ws.Range("O2").Value = 0.2946616623342344
$ws.Range("P2").Value = 0.2946616623342344
$ws.Range("Q2").Value = 0.091124907496
$ws.Range("R2").Value = 0.820124167464
$ws.Range("S2").Value = 0.2946616623342344
$ws.Range("T2").Value = 0.2946616623342344

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.020676
$ws.Range("H3").Value = 0.062028
$ws.Range("O3").Value = 0.2393683991842171
$ws.Range("P3").Value = 0.2393683991842171
$ws.Range("Q3").Value = 0.07402531791999999
$ws.Range("R3").Value = 0.6662278612799999
$ws.Range("S3").Value = 0.2393683991842171
$ws.Range("T3").Value = 0.2393683991842171

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.020676
$ws.Range("H4").Value = 0.062028
$ws.Range("M4").Value = 3.580339
$ws.Range("N4").Value = 10.741017
$ws.Range("O4").Value = 0.2393741266819538
$ws.Range("P4").Value = 0.2393741266819538
$ws.Range("Q4").Value = 0.074027089164
$ws.Range("R4").Value = 0.666243802476
$ws.Range("S4").Value = 0.2393741266819538
$ws.Range("T4").Value = 0.2393741266819538

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.020676
$ws.Range("H5").Value = 0.062028
$ws.Range("M5").Value = 3.389212666666667
$ws.Range("N5").Value = 10.167638
$ws.Range("O5").Value = 0.2265958117995947
$ws.Range("P5").Value = 0.2265958117995947
$ws.Range("Q5").Value = 0.070075361096
$ws.Range("R5").Value = 0.630678249864
$ws.Range("S5").Value = 0.2265958117995947
$ws.Range("T5").Value = 0.2265958117995947
